# Apply updated odds values to row 7 of Sheet1, matching the published
# FlashScore update for the match on row 7 (Atl. Nacional - Envigado).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G7"  = 1.36
    "H7"  = 4.5
    "I7"  = 9
    "J7"  = 1.91
    "L7"  = 8
    "N7"  = 10
    "U7"  = 2.25
    "V7"  = 1.57
    "W7"  = 6
    "X7"  = 6
    "Y7"  = 9
    "Z7"  = 8.5
    "AB7" = 34
    "AC7" = 10
    "AE7" = 23
    "AG7" = 19
    "AJ7" = 101
    "AN7" = 3.2
    "AO7" = 6.5
    "AQ7" = 19
    "AS7" = 151
    "AU7" = 10
    "AV7" = 67
    "AW7" = 9
    "AX7" = 41
    "AY7" = 41
    "AZ7" = 201
    "BA7" = 251
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
